$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2..211 currently hold serial date 45203 (2023-10-04).
# Update them to serial date 45204 (2023-10-05), preserving existing cell formatting/style.
$rng = $ws.Range("C2:C211")
$rng.Value2 = 45204
